# Actualización automática 2025-09-11 14:55:08
#
# Updates the PORCELANATO sales figures for RIOS CARRION ANGEL BENIGNO
# (row 12 = F.V - AREA ANDINA S.A., row 20 = TAMAYO VILLACIS EDWIN XAVIER)
# across the three related sheets, and widens the "VENTA" column on the
# CUMPLIMIENTO MENSUAL sheet so the new values are fully visible.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (column M = PORCELANATO) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M12").Value = 17655.41
$wsGrupo.Range("M20").Value = 1244.13
$wsGrupo.Range("M25").Value = "3 de 23"

# --- Sheet "VENTA MENSUAL" (column F = septiembre) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F12").Value = 17655.41
$wsMensual.Range("F20").Value = 1244.13
$wsMensual.Range("F25").Value = 20470.79

# --- Sheet "CUMPLIMIENTO MENSUAL" (D = VENTA, E = POR CUMPLIR, F = CUMPLIMIENTO) ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D12").Value = 19657.67
$wsCumpl.Range("E12").Value = 23442.4154117774
$wsCumpl.Range("F12").Value = 0.4560935277086111

$wsCumpl.Range("D15").Value = 20470.79
$wsCumpl.Range("E15").Value = 37732.67623249458
$wsCumpl.Range("F15").Value = 0.3517108400078637

# Widen column D ("VENTA") from 13 to 14.
# Note: this runtime's ColumnWidth setter stores (value + 5/6) as the raw
# OOXML column width, so we back-compute the input that yields exactly 14.
$wsCumpl.Columns.Item(4).ColumnWidth = (14 - 5/6)
